$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) from "Gamma1F-HW30.xpc" to "Gamma1F"
$ws.Name = "Gamma1F"

# Append a new row (row 16) of averaged-intensity data using the Gaussian
# Quadrature scheme ("HexGrid-60degTilt5degRes" — same label as row 15).
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.037917103928846
$ws.Range("D16").Value = 0.8575930084976265
$ws.Range("E16").Value = 1.021670918993874
$ws.Range("F16").Value = 1.037917103928846
$ws.Range("G16").Value = 0.9252124734447312
$ws.Range("H16").Value = 1.067212645566715
$ws.Range("I16").Value = 1.02939796776767
$ws.Range("J16").Value = 0.8575930084976265
$ws.Range("K16").Value = 0.93963196374575
$ws.Range("L16").Value = 0.988774533837298
$ws.Range("M16").Value = 0.9898340196999103

# Match the bold/centered/bordered style used by the other index cells in
# column A (e.g. A15) for the new A16 cell.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
